# Weekly fruit/hortaliza price update for "Vega Modelo de Temuco - Zapallo".
# Two new price records are inserted before the existing row 388, shifting
# the remaining 20 records (old rows 388-407) down to rows 390-409.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 388-389 (pushes old 388..407 down to 390..409).
$ws.Rows("388:389").Insert()

# --- New row 388: Camote, 1a nueva(o), origin Perú ---
$ws.Cells.Item(388, 1).Value = 10
$ws.Cells.Item(388, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(388, 3).Value = "La Araucanía"
$ws.Cells.Item(388, 4).Value = 44509
$ws.Cells.Item(388, 5).Value = 9
$ws.Cells.Item(388, 6).Value = 100112045
$ws.Cells.Item(388, 7).Value = "Zapallo"
$ws.Cells.Item(388, 8).Value = "Camote"
$ws.Cells.Item(388, 9).Value = "1a nueva(o)"
$ws.Cells.Item(388, 10).Value = 700
$ws.Cells.Item(388, 11).Value = 800
$ws.Cells.Item(388, 12).Value = 800
$ws.Cells.Item(388, 13).Value = 800
$ws.Cells.Item(388, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(388, 15).Value = "Perú"
$ws.Cells.Item(388, 16).Value = 800
$ws.Cells.Item(388, 17).Value = 1
$ws.Cells.Item(388, 18).Value = "Hortaliza"

# --- New row 389: Paine, 1a (guarda), origin Región del Maule ---
$ws.Cells.Item(389, 1).Value = 10
$ws.Cells.Item(389, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(389, 3).Value = "La Araucanía"
$ws.Cells.Item(389, 4).Value = 44509
$ws.Cells.Item(389, 5).Value = 9
$ws.Cells.Item(389, 6).Value = 100112045
$ws.Cells.Item(389, 7).Value = "Zapallo"
$ws.Cells.Item(389, 8).Value = "Paine"
$ws.Cells.Item(389, 9).Value = "1a (guarda)"
$ws.Cells.Item(389, 10).Value = 1400
$ws.Cells.Item(389, 11).Value = 300
$ws.Cells.Item(389, 12).Value = 400
$ws.Cells.Item(389, 13).Value = 343
$ws.Cells.Item(389, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(389, 15).Value = "Región del Maule"
$ws.Cells.Item(389, 16).Value = 343
$ws.Cells.Item(389, 17).Value = 1
$ws.Cells.Item(389, 18).Value = "Hortaliza"
